$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column price values that are unambiguous text (contain two "." separators)
# -> assign directly, Excel keeps them as text since they are not valid numbers.
$ws.Range("D2").Value = "28.708.51"
$ws.Range("D3").Value = "1.870.02"
$ws.Range("D12").Value = "1.847.52"
$ws.Range("D21").Value = "28.714.38"
$ws.Range("D25").Value = "2.074.17"

# D-column price values that look like plain numbers -> force text format first
# so Excel stores the exact original digits/trailing zeros instead of converting
# them to a floating point number.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "326.70"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4644"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3916"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07901"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9707"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.28"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.741"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.938"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06983"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.39"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.005"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.92"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.51"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.37"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.729"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.001"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "119.45"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09371"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9371"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.322"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.359"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05850"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02126"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.147"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.898"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5659"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.07241"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.76"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5316"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.141"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.137"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.351"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.004"

# E-column percent-change values (always text: leading sign + trailing "%" + padding spaces)
$ws.Range("E2").Value = "  +1.67%  "
$ws.Range("E3").Value = "  +1.74%  "
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("E5").Value = "  -1.16%  "
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("E7").Value = "  +0.81%  "
$ws.Range("E8").Value = "  +1.32%  "
$ws.Range("E9").Value = "  +0.57%  "
$ws.Range("E10").Value = "  +1.08%  "
$ws.Range("E11").Value = "  +1.72%  "
$ws.Range("E12").Value = "  -0.37%  "
$ws.Range("E13").Value = "  +0.48%  "
$ws.Range("E14").Value = "  +0.32%  "
$ws.Range("E15").Value = "  +2.03%  "
$ws.Range("E16").Value = "  +1.70%  "
$ws.Range("E17").Value = "  +0.25%  "
$ws.Range("E18").Value = "  +1.54%  "
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("E20").Value = "  +0.30%  "
$ws.Range("E21").Value = "  +1.56%  "
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("E23").Value = "  +1.17%  "
$ws.Range("E24").Value = "  -0.89%  "
$ws.Range("E25").Value = "  +0.56%  "
$ws.Range("E26").Value = "  +0.22%  "
$ws.Range("E27").Value = "  +1.14%  "
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("E29").Value = "  +0.88%  "
$ws.Range("E30").Value = "  +2.15%  "
$ws.Range("E31").Value = "  +0.94%  "
$ws.Range("E32").Value = "  -0.53%  "
$ws.Range("E33").Value = "  +0.77%  "
$ws.Range("E34").Value = "  +1.65%  "
$ws.Range("E35").Value = "  -2.42%  "
$ws.Range("E36").Value = "  -2.61%  "
$ws.Range("E37").Value = "  -1.27%  "
$ws.Range("E38").Value = "  +0.21%  "
$ws.Range("E39").Value = "  +3.57%  "
$ws.Range("E40").Value = "  +0.67%  "
$ws.Range("E41").Value = "  -0.36%  "
$ws.Range("E42").Value = "  +0.72%  "
$ws.Range("E43").Value = "  +3.12%  "
$ws.Range("E44").Value = "  +1.09%  "
$ws.Range("E45").Value = "  +0.65%  "
$ws.Range("E46").Value = "  -5.50%  "
$ws.Range("E47").Value = "  -8.79%  "
$ws.Range("E48").Value = "  +1.04%  "
$ws.Range("E49").Value = "  +0.85%  "
$ws.Range("E50").Value = "  +0.67%  "
$ws.Range("E51").Value = "  +0.34%  "

Write-Output "cryptos list updated"
